$d = $word.ActiveDocument

# --- Edit 1: "Members Absent: Gerjan Haxhia, ..." -> "... Haxhija, ..." ---
# The target XML shows the sentence split across three runs (identical
# formatting) as if a single letter "j" were typed into the middle of the
# existing word, leaving the text before/after as separate run fragments.
# We insert the "j" as a tracked-changes insertion and then accept it so
# the new run boundary is preserved without any revision markup remaining.
$d.TrackRevisions = $true
$rng = $d.Content
$rng.Find.Execute("Gerjan Haxhi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPoint = $d.Range($rng.End, $rng.End)
$insPoint.InsertAfter("j")
$d.TrackRevisions = $false
$d.Revisions.AcceptAll()

# --- Edit 2: merge "Time: 14:40, " + "April" + " 4" into a single run ---
$d.Content.Find.Execute("Time: 14:40, April 4", $true, $false, $false, $false, $false, $true, 1, $false, "Time: 14:40, April 4", 2)
